$wb = $excel.ActiveWorkbook

$wsNewPlan = $wb.Worksheets.Item("NEW-PLAN")
$wsNewPlan.Range("C2").Value = "PLR310"
$wsNewPlan.Range("E2").Value = "93821778"
$wsNewPlan.Range("D2").Value = "89598071103029413799"
[void]$wsNewPlan.Range("E2").Select()

$wsRes = $wb.Worksheets.Item("POTENTIAL-RESCLIENTS")
$wsRes.Range("B2").Value = "NativoAutomation1"
$wsRes.Range("B3").Value = "NativoAutomation2"
$wsRes.Range("C2").Value = "Prueba"
$wsRes.Range("C3").Value = "Prueba2"
$wsRes.Range("D2").Value = "8877660031"
$wsRes.Range("D3").Value = "8877660032"
[void]$wsRes.Activate()
[void]$wsRes.Range("E7").Select()
